# Vietnamese translation pass for
# "Email 2-1 [TEMPLATE] Partner email - reminder to RSVP.docx"
#
# Each replacement is scoped to the specific paragraph it lives in
# (via $d.Paragraphs(n).Range) so that phrases which repeat verbatim
# elsewhere in the document (e.g. "We hope you're as excited...",
# "We look forward to seeing you at ...", " or ") are not mismatched.
#
# Plain Find.Execute(..., Replace:=...) was observed to (a) make newly
# inserted text inherit the character formatting of an *adjacent*
# w:hyperlink run when the match sits right next to one, and (b) shove a
# structural <w:commentRangeStart/> marker from just-before the replaced
# run to just-after it when the whole run's text is replaced. Both are
# avoided by inserting the new text one character inside the matched
# range (so it lands "inside" the original run, away from any run/marker
# boundary) and then trimming away the leftover original characters.

function SafeReplace($paraIndex, $find, $replace) {
    $d = $word.ActiveDocument
    $p = $d.Paragraphs($paraIndex)
    $scan = $d.Content
    $scan.Start = $p.Range.Start
    $scan.End = $p.Range.End
    $found = $scan.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND: para=$paraIndex find=[$find]"
        return
    }
    $s = $scan.Start
    $e = $scan.End
    $len = $e - $s
    if ($len -lt 2) {
        # Nothing to anchor an interior insert on - direct replace is safe
        # here since these short matches are not next to hyperlinks.
        $scan.Text = $replace
        return
    }
    $rlen = $replace.Length
    $insAt = $d.Range($s + 1, $s + 1)
    $insAt.InsertBefore($replace)
    $trail = $d.Range($s + 1 + $rlen, $e + $rlen)
    $trail.Delete()
    $lead = $d.Range($s, $s + 1)
    $lead.Delete()
}

$d = $word.ActiveDocument

# "Brief" -> "Tóm tắt"   (table cell, paragraph 5)
SafeReplace 5 "Brief" "Tóm tắt"

# "Target audience" -> "Đối tượng mục tiêu"  (table cell, paragraph 8)
SafeReplace 8 "Target audience" "Đối tượng mục tiêu"

# Subject line text (paragraph 12)
SafeReplace 12 ": Reminder: RSVP for " ": Nhắc nhở: Xác nhận tham dự sự kiện "

# "Don't delay! Book your spot today!" (paragraph 15)
SafeReplace 15 "Don’t delay! Book your spot today!" "Đừng chậm trễ! Hãy đặt chỗ của bạn ngay hôm nay!"

# "Hi " -> "Xin chào " (paragraph 17)
SafeReplace 17 "Hi " "Xin chào "

# Paragraph 19: "We hope you're as excited as us for the [EVENT NAME], happening on [DD Mmm YYYY]!"
SafeReplace 19 "We hope you’re as excited as us for the " "Chúng tôi hy vọng bạn mong chờ sự kiện "
SafeReplace 19 ", happening on " " diễn ra vào ngày "
SafeReplace 19 "!" " như chúng tôi!"

# Paragraph 20: "We hope you're as excited as us for the [EVENT NAME], happening from [DD Mmm YYYY] to [DD Mmm YYYY]!"
SafeReplace 20 "We hope you’re as excited as us for the " "Chúng tôi hy vọng bạn mong chờ sự kiện "
SafeReplace 20 ", happening from " " diễn ra từ ngày "
SafeReplace 20 " to " " đến ngày "
SafeReplace 20 "!" " như chúng tôi!"

# Paragraph 21: "Confirm your attendance for this highly-anticipated event by [DD Mmm YYYY] as spots are limited and on a first-come, first-served basis."
SafeReplace 21 "Confirm your attendance for this highly-anticipated event by [" "Hãy xác nhận tham dự sự kiện rất được mong đợi này trước ngày ["
SafeReplace 21 "] as spots are limited and on a first-come, first-served basis." "] vì số lượng tham dự là có hạn và chúng tôi sẽ chọn những đơn đăng ký sớm nhất."

# "RVSP now" button text (paragraph 22)
SafeReplace 22 "RVSP now" "Xác nhận tham dự ngay bây giờ"

# Paragraph 25: "If you have any questions, please contact us via live chat or WhatsApp. "
SafeReplace 25 "If you have any questions, please contact us via " "Nếu bạn cần hỗ trợ, hãy liên hệ với chúng tôi qua "
SafeReplace 25 " or " " hoặc "

# Paragraph 26: "If you have any questions, please contact your country manager, [NAME], at [EMAIL ADDRESS] or [WHATSAPP NO] (WhatsApp). "
SafeReplace 26 "If you have any questions, please contact your country manager, " "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn "
SafeReplace 26 ", at " ", qua email "
SafeReplace 26 " or " " hoặc số "

# Paragraph 27: "We look forward to seeing you at [EVENT NAME]! " ([EVENT NAME] is its own run)
SafeReplace 27 "We look forward to seeing you at " "Chúng tôi rất mong được gặp bạn tại sự kiện "

# Paragraph 38: "We look forward to seeing you at [EVENT NAME]! " (literal text, single run)
SafeReplace 38 "We look forward to seeing you at [EVENT NAME]! " "Chúng tôi rất mong được gặp bạn tại sự kiện [EVENT NAME]! "

# Paragraph 39: "If you have any questions, please contact your country manager:"
SafeReplace 39 "If you have any questions, please contact your country manager:" "Nếu bạn có bất kỳ thắc mắc nào, vui lòng liên hệ với giám đốc phụ trách quốc gia của bạn:"

# Paragraph 41: "If you have any questions, please contact us via:"
SafeReplace 41 "If you have any questions, please contact us via:" "Nếu bạn cần hỗ trợ, vui lòng liên hệ với chúng tôi qua:"

Write-Host "Translation replacements applied."
